$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.213.20'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '1.600.41'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = "'303.41"
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = "'52.04"
$ws.Range("E8").Value = '  +4.92%  '
$ws.Range("D9").Value = "'0.3633"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = "'0.08136"
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").Value = "'22.70"
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = "'6.568"
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = "'7.399"
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("D17").Value = '1.601.43'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = "'94.27"
$ws.Range("E18").Value = '  +2.47%  '
$ws.Range("D19").Value = "'0.06927"
$ws.Range("E19").Value = '  +1.96%  '
$ws.Range("D20").Value = "'18.08"
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = "'6.519"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = "'12.88"
$ws.Range("E23").Value = '  -1.65%  '
$ws.Range("D24").Value = '23.219.28'
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").Value = "'2.422"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("D26").Value = "'3.027"
$ws.Range("E26").Value = '  +7.06%  '
$ws.Range("D27").Value = "'21.17"
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").Value = "'149.33"
$ws.Range("D29").Value = "'5.280"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("D30").Value = "'134.97"
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").Value = "'2.389"
$ws.Range("E31").Value = '  +6.79%  '
$ws.Range("D32").Value = "'6.700"
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("D33").Value = '1.778.00'
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").Value = "'0.9634"
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("D35").Value = "'0.07466"
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("D36").Value = "'10.36"
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = "'0.02738"
$ws.Range("E37").Value = '  +1.69%  '
$ws.Range("D38").Value = "'0.2529"
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").Value = "'0.08788"
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("D40").Value = "'6.085"
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("D42").Value = "'0.7074"
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").Value = "'12.37"
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = "'15.55"
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").Value = "'0.6519"
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = "'2.312"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("D49").Value = "'132.29"
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").Value = "'0.07916"
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").Value = "'1.203"
$ws.Range("E51").Value = '  -0.25%  '
